$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update description texts: merge the two-line descriptions into a single
#     line (replace embedded newline with a space) for Magicien, Eveque,
#     Marchand and Architecte rows ---
$ws.Range("D4").Value = "Au choix: Le ¨Magicien peut échanger la totalité de ses cartes avec le joueur de son choix. Ou: le Magicien peut échanger des cartes de sa main contre le même nombre de cartes de la pioche."
$ws.Range("D6").Value = "L'Évêque ne peut pas être attaqué par le Condottière. Chaque quartier religieux qu'il possède lui rapporte une pièce d'or."
$ws.Range("D7").Value = "Le Marchand reçoit une pièce d'or en plus au début de son tour. Chaque quartier marchand qu'il possède lui rapporte une pièce d'or."
$ws.Range("D8").Value = "L'Architecte pioche deux cartes quartier en plus. il peut bâtir jusqu'à trois quartiers."

# --- Apply base alignment (left / top) to the whole used range ---
$all = $ws.Range("A1:D9")
$all.HorizontalAlignment = -4131   # xlLeft
$all.VerticalAlignment = -4160     # xlTop

# --- Description column: wrap the long text cells (all except the header
#     and the Voleur row, which stays unwrapped). Multi-area ("D2,D4:D9")
#     ranges are not reliable, so the two contiguous pieces are handled
#     separately. ---
$ws.Range("D2").WrapText = $true
$ws.Range("D4:D9").WrapText = $true

# --- Row heights: only Magicien (row 4) and Condottiere (row 9) keep the
#     taller 30pt row; the others go back to the default height ---
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 30

# --- Update the active selection to D6 ---
$ws.Range("D6").Select()

$wb.Save()
